$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Good Morning" greeting cell for rule R10 to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the active/selected cell, matching the saved selection state
$ws.Range("E8").Select()
